# Update the dSF column (F) values for rows 2-7 as per repulled data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -3
$ws.Range("F3").Value = -6
$ws.Range("F4").Value = 2
$ws.Range("F5").Value = 2
$ws.Range("F6").Value = 3
$ws.Range("F7").Value = -1
